$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.466.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4819"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2803"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.856.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07449"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.095"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6421"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.449.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007487"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.100.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.149"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.109"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.354"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.907"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1056"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +13.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.386"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.39%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.991"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04985"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.183"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7440"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01942"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.639"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9170"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4200"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.576"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.231"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1229"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.901"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.425"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
